# ProjectApplication.xlsx — "A new massive change / Some MVC done"
#
# 1. The applicant's NRIC on the first data row (row 2) was corrected and the
#    submission timestamp bumped.
# 2. Two new blank rows were inserted above the last record (which pushes the
#    old row 4 down to row 6).
# 3. The (now) last record's Project ID, Applicant NRIC and timestamp were
#    updated to reflect a new application.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Applicant NRIC corrected, Date refreshed -----------------------
$ws.Range("C2").Value = "T2109876H"
$ws.Range("F2").Value = 45769.835042048609

# --- Make room for the two new blank rows (old row 4 becomes row 6) --------
$ws.Rows("4:5").Insert()

# --- Row 6 (previously row 4): new applicant record -------------------------
$ws.Range("B6").Value = 4
$ws.Range("C6").Value = "T2345678D"
$ws.Range("F6").Value = 45769.900285069445

# --- Match the workbook's on-screen selection at save time ------------------
$ws.Range("A4:F5").Select()

$wb.Save()
